$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5273.2666
$ws.Range("I32").Value = 4685.5557
$ws.Range("J32").Value = 6154.8335
$ws.Range("K32").Value = 4685.5557
$ws.Range("L32").Value = 6154.8335
$ws.Range("M32").Value = -4359.5557
$ws.Range("N32").Value = -6806.8335

$ws.Range("H41").Value = 873
$ws.Range("I41").Value = 935.25
$ws.Range("J41").Value = 748.5
$ws.Range("K41").Value = 935.25
$ws.Range("L41").Value = 748.5
$ws.Range("M41").Value = -495.25

$ws.Range("H51").Value = 7133.1665
$ws.Range("I51").Value = 3199.75
$ws.Range("J51").Value = 15000
$ws.Range("K51").Value = 3199.75
$ws.Range("L51").Value = 15000
$ws.Range("M51").Value = -2715.75

$ws.Range("H55").Value = 464.9
$ws.Range("I55").Value = 310.25
$ws.Range("J55").Value = 568
$ws.Range("K55").Value = 310.25
$ws.Range("L55").Value = 568
$ws.Range("M55").Value = -96.25

$ws.Range("H70").Value = 3960.5
$ws.Range("I70").Value = 7550.5
$ws.Range("J70").Value = 1567.1666
$ws.Range("K70").Value = 22651.5
$ws.Range("L70").Value = 4701.4998
$ws.Range("M70").Value = -22381.5
$ws.Range("N70").Value = -5241.4998

$ws.Range("H73").Value = 3960.5
$ws.Range("I73").Value = 7550.5
$ws.Range("J73").Value = 1567.1666
$ws.Range("K73").Value = 22651.5
$ws.Range("L73").Value = 4701.4998
$ws.Range("M73").Value = -21715.5
$ws.Range("N73").Value = -6573.4998

$ws.Range("H86").Value = 5849.5
$ws.Range("I86").Value = 5800
$ws.Range("J86").Value = 5899
$ws.Range("K86").Value = 5800
$ws.Range("L86").Value = 5899
$ws.Range("M86").Value = -4677
$ws.Range("N86").Value = -8145

$ws.Range("H89").Value = 5849.5
$ws.Range("I89").Value = 5800
$ws.Range("J89").Value = 5899
$ws.Range("K89").Value = 29000
$ws.Range("L89").Value = 29495
$ws.Range("M89").Value = -23384
$ws.Range("N89").Value = -40727

$ws.Range("H100").Value = 2963.3635
$ws.Range("I100").Value = 3200.1428
$ws.Range("J100").Value = 2549
$ws.Range("K100").Value = 3200.1428
$ws.Range("L100").Value = 2549
$ws.Range("M100").Value = -2659.1428
$ws.Range("N100").Value = -3631

$ws.Range("H113").Value = 9385.286
$ws.Range("I113").Value = 7833
$ws.Range("J113").Value = 10549.5
$ws.Range("K113").Value = 7833
$ws.Range("L113").Value = 10549.5
$ws.Range("M113").Value = -4579

$ws.Range("H121").Value = 400
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 400
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 1200
$ws.Range("N121").Value = -4694

$ws.Range("H137").Value = 2932.111
$ws.Range("I137").Value = 1607.1666
$ws.Range("J137").Value = 5582
$ws.Range("K137").Value = 4821.4998
$ws.Range("L137").Value = 16746
$ws.Range("M137").Value = -2271.4998

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1404
$ws.Range("I74").Value = 1404
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1404
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -530
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1404
$ws.Range("I77").Value = 1404
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 7020
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -2652
$ws.Range("N77").ClearContents()

$ws.Range("H110").Value = 4355.8667
$ws.Range("I110").Value = 4355.8667
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 4355.8667
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = -2310.8667

$ws.Range("H122").Value = 2829.5
$ws.Range("I122").Value = 2829.5
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 8488.5
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -6038.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 4748
$ws.Range("I20").Value = 4250
$ws.Range("J20").Value = 4997
$ws.Range("K20").Value = 4250
$ws.Range("L20").Value = 4997
$ws.Range("M20").Value = -4003
$ws.Range("N20").Value = -5491

$ws.Range("H86").Value = 3052.5557
$ws.Range("I86").Value = 3157.8667
$ws.Range("J86").Value = 2526
$ws.Range("K86").Value = 3157.8667
$ws.Range("L86").Value = 2526
$ws.Range("M86").Value = -2034.8667

$ws.Range("H89").Value = 3052.5557
$ws.Range("I89").Value = 3157.8667
$ws.Range("J89").Value = 2526
$ws.Range("K89").Value = 15789.3335
$ws.Range("L89").Value = 12630
$ws.Range("M89").Value = -10173.3335

$ws.Range("H105").Value = 5746.6113
$ws.Range("I105").Value = 5574.2856
$ws.Range("J105").Value = 5856.273
$ws.Range("K105").Value = 5574.2856
$ws.Range("L105").Value = 5856.273
$ws.Range("M105").Value = -3827.2856
$ws.Range("N105").Value = -9350.273000000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H47").Value = 19100
$ws.Range("I47").Value = 19100
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 19100
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -18534

$ws.Range("H62").Value = 2997.5
$ws.Range("I62").Value = 2997
$ws.Range("J62").Value = 2999
$ws.Range("K62").Value = 2997
$ws.Range("L62").Value = 2999
$ws.Range("M62").Value = -2373

$ws.Range("H65").Value = 2997.5
$ws.Range("I65").Value = 2997
$ws.Range("J65").Value = 2999
$ws.Range("K65").Value = 14985
$ws.Range("L65").Value = 14995
$ws.Range("M65").Value = -11865

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1199.6666
$ws.Range("I5").Value = 0
$ws.Range("J5").Value = 1199.6666
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 3598.9998
$ws.Range("N5").Value = -3822.9998

$ws.Range("H29").Value = 200.75
$ws.Range("I29").Value = 225.5
$ws.Range("J29").Value = 176
$ws.Range("K29").Value = 676.5
$ws.Range("L29").Value = 528
$ws.Range("M29").Value = -399.5
$ws.Range("N29").Value = -1082

$ws.Range("H135").Value = 1199.6666
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 1199.6666
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 10796.9994
$ws.Range("N135").Value = -15866.9994

$ws.Range("H136").Value = 8749.75
$ws.Range("I136").Value = 8749.75
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 26249.25
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -21149.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6041.4287
$ws.Range("I70").Value = 6728.2
$ws.Range("J70").Value = 4324.5
$ws.Range("K70").Value = 6728.2
$ws.Range("L70").Value = 4324.5
$ws.Range("M70").Value = -6458.2

$ws.Range("H73").Value = 6041.4287
$ws.Range("I73").Value = 6728.2
$ws.Range("J73").Value = 4324.5
$ws.Range("K73").Value = 6728.2
$ws.Range("L73").Value = 4324.5
$ws.Range("M73").Value = -5792.2

$ws.Range("H80").Value = 2636.875
$ws.Range("I80").Value = 2875
$ws.Range("J80").Value = 2398.75
$ws.Range("K80").Value = 2875
$ws.Range("L80").Value = 2398.75
$ws.Range("M80").Value = -1877
$ws.Range("N80").Value = -4394.75

$ws.Range("H83").Value = 2636.875
$ws.Range("I83").Value = 2875
$ws.Range("J83").Value = 2398.75
$ws.Range("K83").Value = 14375
$ws.Range("L83").Value = 11993.75
$ws.Range("M83").Value = -9383
$ws.Range("N83").Value = -21977.75

$ws.Range("H97").Value = 550.7692
$ws.Range("I97").Value = 522.2222
$ws.Range("J97").Value = 615
$ws.Range("K97").Value = 522.2222
$ws.Range("L97").Value = 615
$ws.Range("M97").Value = -26.22220000000004
$ws.Range("N97").Value = -1607

$ws.Range("H122").Value = 4303.4287
$ws.Range("I122").Value = 3500.3333
$ws.Range("J122").Value = 4905.75
$ws.Range("K122").Value = 10500.9999
$ws.Range("L122").Value = 14717.25
$ws.Range("M122").Value = -8050.999899999999

$ws.Range("H128").Value = 44997.5
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 44997.5
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 44997.5
$ws.Range("N128").Value = -54957.5
$ws.Range("M128").ClearContents()

$ws.Range("H132").Value = 2000
$ws.Range("I132").Value = 1999.6666
$ws.Range("J132").Value = 2001
$ws.Range("K132").Value = 5998.9998
$ws.Range("L132").Value = 6003
$ws.Range("M132").Value = -3468.9998
$ws.Range("N132").Value = -11063

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 769.125
$ws.Range("I16").Value = 607.7143
$ws.Range("J16").Value = 1899
$ws.Range("K16").Value = 607.7143
$ws.Range("L16").Value = 1899
$ws.Range("M16").Value = -437.7143
$ws.Range("N16").Value = -2239

$ws.Range("H22").Value = 1381.9166
$ws.Range("I22").Value = 875.8889
$ws.Range("J22").Value = 2900
$ws.Range("K22").Value = 875.8889
$ws.Range("L22").Value = 2900
$ws.Range("M22").Value = -580.8889

$ws.Range("H27").Value = 1381.9166
$ws.Range("I27").Value = 875.8889
$ws.Range("J27").Value = 2900
$ws.Range("K27").Value = 875.8889
$ws.Range("L27").Value = 2900
$ws.Range("M27").Value = -768.8889

$ws.Range("H30").Value = 1027.6666
$ws.Range("I30").Value = 1163.2
$ws.Range("J30").Value = 350
$ws.Range("K30").Value = 1163.2
$ws.Range("L30").Value = 350
$ws.Range("M30").Value = -1055.2

$ws.Range("H55").Value = 807.4286
$ws.Range("I55").Value = 613.125
$ws.Range("J55").Value = 1066.5
$ws.Range("K55").Value = 613.125
$ws.Range("L55").Value = 1066.5
$ws.Range("M55").Value = -440.125
$ws.Range("N55").Value = -1412.5

$ws.Range("H75").Value = 36000
$ws.Range("I75").Value = 36000
$ws.Range("J75").Value = 0
$ws.Range("K75").Value = 36000
$ws.Range("L75").Value = 0
$ws.Range("M75").Value = -35064

$ws.Range("H78").Value = 36000
$ws.Range("I78").Value = 36000
$ws.Range("J78").Value = 0
$ws.Range("K78").Value = 108000
$ws.Range("L78").Value = 0
$ws.Range("M78").Value = -103320

$ws.Range("H82").Value = 3133.3333
$ws.Range("I82").Value = 2700
$ws.Range("J82").Value = 4000
$ws.Range("K82").Value = 2700
$ws.Range("L82").Value = 4000
$ws.Range("M82").Value = -2339
$ws.Range("N82").Value = -4722

$ws.Range("H85").Value = 3133.3333
$ws.Range("I85").Value = 2700
$ws.Range("J85").Value = 4000
$ws.Range("K85").Value = 2700
$ws.Range("L85").Value = 4000
$ws.Range("M85").Value = -1452
$ws.Range("N85").Value = -6496

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H45").Value = 40109.332
$ws.Range("I45").Value = 31198
$ws.Range("J45").Value = 51248.5
$ws.Range("K45").Value = 31198
$ws.Range("L45").Value = 51248.5
$ws.Range("M45").Value = -30707
$ws.Range("N45").Value = -52230.5

$ws.Range("H81").Value = 3621.7778
$ws.Range("I81").Value = 2156.5715
$ws.Range("J81").Value = 8750
$ws.Range("K81").Value = 4313.143
$ws.Range("L81").Value = 17500
$ws.Range("M81").Value = -3252.143
$ws.Range("N81").Value = -19622

$ws.Range("H84").Value = 3621.7778
$ws.Range("I84").Value = 2156.5715
$ws.Range("J84").Value = 8750
$ws.Range("K84").Value = 21565.715
$ws.Range("L84").Value = 87500
$ws.Range("M84").Value = -16261.715
$ws.Range("N84").Value = -98108

$ws.Range("H122").Value = 2624.7368
$ws.Range("I122").Value = 2578.6667
$ws.Range("J122").Value = 2797.5
$ws.Range("K122").Value = 7736.000100000001
$ws.Range("L122").Value = 8392.5
$ws.Range("M122").Value = -5286.000100000001
$ws.Range("N122").Value = -13292.5

$ws.Range("H126").Value = 2144.5557
$ws.Range("I126").Value = 2056.5
$ws.Range("J126").Value = 2849
$ws.Range("K126").Value = 6169.5
$ws.Range("L126").Value = 8547
$ws.Range("M126").Value = -3699.5

$ws.Range("H136").Value = 3922.4102
$ws.Range("I136").Value = 4156.6665
$ws.Range("J136").Value = 3649.111
$ws.Range("K136").Value = 12469.9995
$ws.Range("L136").Value = 10947.333
$ws.Range("M136").Value = -9919.999500000002
$ws.Range("N136").Value = -16047.333
